$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet
$ws.Range("A103").Font.Name = "Calibri"
$ws.Range("A103").Font.Size = 11
$ws.Range("A103").NumberFormat = "yyyy-mm-dd"
